# ---------------------------------------------------------------------------
# karma_performance.xlsx — benchmark and doc updates
#
# 1. Rename the three sheets and move the active tab to "Single int".
# 2. Split the data that used to live twice in Sheet1 (rows 3-10 and
#    36-42) so that each benchmark gets its own sheet:
#      Sheet1 "Single double"      -> keeps rows 3-10
#      Sheet2 "Sequence of items"  -> gets the old rows 36-42 (as rows 3-9)
#      Sheet3 "Single int"         -> brand-new benchmark data
# 3. Re-point the two existing charts (chart1/chart2) at their new sheet
#    names / ranges, and build a third chart for the "Single int" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. rename sheets -------------------------------------------------
$wsSingleDouble = $wb.Worksheets.Item("Sheet1")
$wsSeqItems     = $wb.Worksheets.Item("Sheet2")
$wsSingleInt    = $wb.Worksheets.Item("Sheet3")

$wsSingleDouble.Name = "Single double"
$wsSeqItems.Name     = "Sequence of items"
$wsSingleInt.Name    = "Single int"

# --- 2. move the "sequence of items" block from Sheet1 to Sheet2 ------
$srcBlock = $wsSingleDouble.Range("D36:P42")
$srcBlock.Copy()
$wsSeqItems.Range("D3").PasteSpecial(-4104)
$wsSingleDouble.Application.CutCopyMode = $false
$srcBlock.Clear()

$wsSeqItems.Range("K1:K1048576").Select()

# --- 3. add the "Single int" benchmark data ----------------------------
$wsSingleInt.Range("D1").Value = "Converting 10000000 randomly generated int values to strings."

$wsSingleInt.Range("E3").Value = "gcc 4.4.0 (32)"
$wsSingleInt.Range("F3").Value = "VC++ 10 (32)"
$wsSingleInt.Range("G3").Value = "Intel 11.1 (32)"
$wsSingleInt.Range("H3").Value = "gcc 4.4.0 (64)"
$wsSingleInt.Range("I3").Value = "VC++ 10 (64)"
$wsSingleInt.Range("J3").Value = "Intel 11.1 (64)"

$wsSingleInt.Range("D4").Value = "ltoa "
$wsSingleInt.Range("E4").Value = 1.5589999999999999
$wsSingleInt.Range("F4").Value = 0.88100000000000001
$wsSingleInt.Range("G4").Value = 0.85899999999999999
$wsSingleInt.Range("H4").Value = 1.2
$wsSingleInt.Range("I4").Value = 1.1180000000000001
$wsSingleInt.Range("J4").Value = 0.88900000000000001

$wsSingleInt.Range("D5").Value = "iostreams   "
$wsSingleInt.Range("E5").Value = 6.484
$wsSingleInt.Range("F5").Value = 13.161
$wsSingleInt.Range("G5").Value = 11.635999999999999
$wsSingleInt.Range("H5").Value = 3.42
$wsSingleInt.Range("I5").Value = 7.8120000000000003
$wsSingleInt.Range("J5").Value = 7.3680000000000003

$wsSingleInt.Range("D6").Value = "Boost.Format"
$wsSingleInt.Range("E6").Value = 16.823
$wsSingleInt.Range("F6").Value = 21.568999999999999
$wsSingleInt.Range("G6").Value = 19.706
$wsSingleInt.Range("H6").Value = 17.28
$wsSingleInt.Range("I6").Value = 14.401999999999999
$wsSingleInt.Range("J6").Value = 13.222

$wsSingleInt.Range("D7").Value = "Karma"
$wsSingleInt.Range("E7").Value = 2.5619999999999998
$wsSingleInt.Range("F7").Value = 1.0109999999999999
$wsSingleInt.Range("G7").Value = 0.95499999999999996
$wsSingleInt.Range("H7").Value = 2.956
$wsSingleInt.Range("I7").Value = 1.016
$wsSingleInt.Range("J7").Value = 0.878

# column widths matching the other two data sheets
$wsSingleInt.Columns.Item(4).ColumnWidth = 13.5703125
$wsSingleInt.Columns.Item(5).ColumnWidth = 12.140625
$wsSingleInt.Columns.Item(6).ColumnWidth = 11.7109375
$wsSingleInt.Columns.Item(7).ColumnWidth = 13.140625
$wsSingleInt.Columns.Item(8).ColumnWidth = 12.140625
$wsSingleInt.Columns.Item(9).ColumnWidth = 11.7109375
$wsSingleInt.Columns.Item(10).ColumnWidth = 13.140625

# wrapped, otherwise-empty cell below the data (matches styles.xml xf#1)
$wsSingleInt.Range("E12").WrapText = $true
$wsSingleInt.Range("P7:P10").Select()

# --- 4. re-point chart1 ("Single double") at the renamed sheet --------
$chart1 = $wsSingleDouble.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = "=SERIES('Single double'!`$D`$4,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$4:`$J`$4,1)"
$chart1.SeriesCollection().Item(2).Formula = "=SERIES('Single double'!`$D`$5,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$5:`$J`$5,2)"
$chart1.SeriesCollection().Item(3).Formula = "=SERIES('Single double'!`$D`$6,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$6:`$J`$6,3)"
$chart1.SeriesCollection().Item(4).Formula = "=SERIES('Single double'!`$D`$7,'Single double'!`$E`$3:`$J`$3,'Single double'!`$E`$7:`$J`$7,4)"

# --- 5. move/re-point chart2 ("Sequence of items") --------------------
# chart2 used to live as the second chart object on Sheet1; rebuild the
# same four series directly on the "Sequence of items" sheet and drop
# the old one.
$oldChart2Obj = $wsSingleDouble.ChartObjects().Item(1)
# (chart1 is Item(1) after nothing was removed yet on Sheet1 -- the
#  "sequence of items" chart is the remaining second chart object)
$seqChartObj = $wsSingleDouble.ChartObjects().Item(2)
$chart2 = $seqChartObj.Chart

$chart2.SeriesCollection().Item(1).Formula = "=SERIES('Sequence of items'!`$D`$4,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$4:`$J`$4,1)"
$chart2.SeriesCollection().Item(2).Formula = "=SERIES('Sequence of items'!`$D`$5,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$5:`$J`$5,2)"
$chart2.SeriesCollection().Item(3).Formula = "=SERIES('Sequence of items'!`$D`$6,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$6:`$J`$6,3)"
$chart2.SeriesCollection().Item(4).Formula = "=SERIES('Sequence of items'!`$D`$7,'Sequence of items'!`$E`$3:`$J`$3,'Sequence of items'!`$E`$7:`$J`$7,4)"

$seqChartObj.Copy()
$wsSeqItems.Range("E12").Select()
$wsSeqItems.Paste()
$oldChart2Obj = $null
$seqChartObj.Delete()

# --- 6. build chart3 ("Single int") ------------------------------------
$intChartObj = $wsSingleInt.ChartObjects().Add(66675, 44450, 2857500, 1724025)
$chart3 = $intChartObj.Chart
$chart3.ChartType = 51

$s1 = $chart3.SeriesCollection().NewSeries()
$s1.Formula = "=SERIES('Single int'!`$D`$4,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$4:`$J`$4,1)"
$s2 = $chart3.SeriesCollection().NewSeries()
$s2.Formula = "=SERIES('Single int'!`$D`$5,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$5:`$J`$5,2)"
$s3 = $chart3.SeriesCollection().NewSeries()
$s3.Formula = "=SERIES('Single int'!`$D`$6,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$6:`$J`$6,3)"
$s4 = $chart3.SeriesCollection().NewSeries()
$s4.Formula = "=SERIES('Single int'!`$D`$7,'Single int'!`$E`$3:`$J`$3,'Single int'!`$E`$7:`$J`$7,4)"

$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Format single int" + [char]10 + "(10000000 iterations)"

$chart3.ChartGroups(1).GapWidth = 152

$catAx3 = $chart3.Axes(1)
$catAx3.HasMajorGridlines = $true

$valAx3 = $chart3.Axes(2)
$valAx3.HasMajorGridlines = $true
$valAx3.HasTitle = $true
$valAx3.AxisTitle.Text = "Measured time [s]"

$chart3.HasLegend = $true
$chart3.Legend.Position = -4107

# --- 7. selections / active sheet --------------------------------------
$wsSingleDouble.Range("D3").Select()
$wsSingleInt.Activate()
